# Auto-generated edit script applying cached-value updates to the
# "Leve Profits" calculation sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# These are static numeric snapshots (no formulas in the sheets), refreshed
# by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 698.61536
$ws.Range("I28").Value = 698.61536
$ws.Range("K28").Value = 698.61536
$ws.Range("M28").Value = -213.61536
$ws.Range("H53").Value = 1169.12
$ws.Range("I53").Value = 951.2857
$ws.Range("J53").Value = 1446.3636
$ws.Range("K53").Value = 951.2857
$ws.Range("L53").Value = 1446.3636
$ws.Range("M53").Value = -314.2857
$ws.Range("N53").Value = -2720.3636
$ws.Range("H58").Value = 2925.8
$ws.Range("I58").Value = 3600
$ws.Range("J58").Value = 1914.5
$ws.Range("K58").Value = 10800
$ws.Range("L58").Value = 5743.5
$ws.Range("M58").Value = -10650
$ws.Range("N58").Value = -6043.5
$ws.Range("H132").Value = 2856.9153
$ws.Range("I132").Value = 2814.9216
$ws.Range("J132").Value = 3124.625
$ws.Range("K132").Value = 8444.764800000001
$ws.Range("L132").Value = 9373.875
$ws.Range("M132").Value = -5914.764800000001
$ws.Range("N132").Value = -14433.875

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10419559
$ws.Range("I32").Value = 5556912.5
$ws.Range("J32").Value = 25007500
$ws.Range("K32").Value = 5556912.5
$ws.Range("L32").Value = 25007500
$ws.Range("M32").Value = -5556625.5
$ws.Range("N32").Value = -25008074
$ws.Range("H45").Value = 4371.05
$ws.Range("I45").Value = 3856.7778
$ws.Range("J45").Value = 8999.5
$ws.Range("K45").Value = 3856.7778
$ws.Range("L45").Value = 8999.5
$ws.Range("M45").Value = -3479.7778
$ws.Range("N45").Value = -9753.5
$ws.Range("H63").Value = 3569.8333
$ws.Range("I63").Value = 2255.4285
$ws.Range("K63").Value = 2255.4285
$ws.Range("M63").Value = -1569.4285
$ws.Range("H66").Value = 3569.8333
$ws.Range("I66").Value = 2255.4285
$ws.Range("K66").Value = 11277.1425
$ws.Range("M66").Value = -7845.1425
$ws.Range("H74").Value = 3999.6667
$ws.Range("J74").Value = 5999.5
$ws.Range("L74").Value = 5999.5
$ws.Range("N74").Value = -7747.5
$ws.Range("H77").Value = 3999.6667
$ws.Range("J77").Value = 5999.5
$ws.Range("L77").Value = 29997.5
$ws.Range("N77").Value = -38733.5
$ws.Range("H97").Value = 1029.963
$ws.Range("I97").Value = 692.2273
$ws.Range("J97").Value = 2516
$ws.Range("K97").Value = 692.2273
$ws.Range("L97").Value = 2516
$ws.Range("M97").Value = -196.2273
$ws.Range("N97").Value = -3508
$ws.Range("H102").Value = 2535
$ws.Range("I102").Value = 2145.6
$ws.Range("K102").Value = 2145.6
$ws.Range("M102").Value = -523.5999999999999
$ws.Range("H132").Value = 3286.5833
$ws.Range("I132").Value = 3030.56
$ws.Range("K132").Value = 9091.68
$ws.Range("M132").Value = -6561.68

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2325.762
$ws.Range("I86").Value = 2227.889
$ws.Range("K86").Value = 2227.889
$ws.Range("M86").Value = -1104.889
$ws.Range("H89").Value = 2325.762
$ws.Range("I89").Value = 2227.889
$ws.Range("K89").Value = 11139.445
$ws.Range("M89").Value = -5523.445
$ws.Range("H94").Value = 997.86957
$ws.Range("I94").Value = 534.75
$ws.Range("J94").Value = 2056.4285
$ws.Range("K94").Value = 534.75
$ws.Range("L94").Value = 2056.4285
$ws.Range("M94").Value = -83.75
$ws.Range("N94").Value = -2958.4285
$ws.Range("H105").Value = 3456.7646
$ws.Range("I105").Value = 2877.375
$ws.Range("J105").Value = 3971.7778
$ws.Range("K105").Value = 2877.375
$ws.Range("L105").Value = 3971.7778
$ws.Range("M105").Value = -1130.375
$ws.Range("N105").Value = -7465.7778
$ws.Range("H107").Value = 1804.875
$ws.Range("I107").Value = 1562.7142
$ws.Range("K107").Value = 1562.7142
$ws.Range("M107").Value = 357.2858000000001
$ws.Range("H135").Value = 60000
$ws.Range("J135").Value = 60000
$ws.Range("L135").Value = 60000
$ws.Range("N135").Value = -70140

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4964.619
$ws.Range("J31").Value = 5990.5
$ws.Range("L31").Value = 5990.5
$ws.Range("N31").Value = -6580.5
$ws.Range("H34").Value = 4964.619
$ws.Range("J34").Value = 5990.5
$ws.Range("L34").Value = 5990.5
$ws.Range("N34").Value = -6394.5
$ws.Range("H58").Value = 2430.5881
$ws.Range("I58").Value = 1986.8148
$ws.Range("K58").Value = 1986.8148
$ws.Range("M58").Value = -1783.8148
$ws.Range("H62").Value = 1990
$ws.Range("I62").Value = 1990
$ws.Range("K62").Value = 1990
$ws.Range("M62").Value = -1366
$ws.Range("H65").Value = 1990
$ws.Range("I65").Value = 1990
$ws.Range("K65").Value = 9950
$ws.Range("M65").Value = -6830
$ws.Range("H99").Value = 1717.25
$ws.Range("I99").Value = 1717.25
$ws.Range("K99").Value = 1717.25
$ws.Range("M99").Value = -219.25
$ws.Range("H105").Value = 1906.125
$ws.Range("I105").Value = 1986.5
$ws.Range("J105").Value = 1665
$ws.Range("K105").Value = 1986.5
$ws.Range("L105").Value = 1665
$ws.Range("M105").Value = -239.5
$ws.Range("N105").Value = -5159
$ws.Range("H107").Value = 79318.766
$ws.Range("I107").Value = 144203
$ws.Range("K107").Value = 144203
$ws.Range("M107").Value = -142283
$ws.Range("H126").Value = 1717.25
$ws.Range("I126").Value = 1717.25
$ws.Range("K126").Value = 5151.75
$ws.Range("M126").Value = -2681.75
$ws.Range("H132").Value = 1813.875
$ws.Range("I132").Value = 751.8333
$ws.Range("K132").Value = 2255.4999
$ws.Range("M132").Value = 274.5001000000002
$ws.Range("H134").Value = 2007.973
$ws.Range("I134").Value = 1129.9286
$ws.Range("J134").Value = 4739.6665
$ws.Range("K134").Value = 3389.7858
$ws.Range("L134").Value = 14218.9995
$ws.Range("M134").Value = -854.7857999999997
$ws.Range("N134").Value = -19288.9995
$ws.Range("H136").Value = 2430.5881
$ws.Range("I136").Value = 1986.8148
$ws.Range("K136").Value = 5960.4444
$ws.Range("M136").Value = -3410.4444

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 569.375
$ws.Range("J38").Value = 874.7
$ws.Range("L38").Value = 2624.1
$ws.Range("N38").Value = -3318.1

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4416.6665
$ws.Range("I70").Value = 3500
$ws.Range("K70").Value = 3500
$ws.Range("M70").Value = -3230
$ws.Range("H73").Value = 4416.6665
$ws.Range("I73").Value = 3500
$ws.Range("K73").Value = 3500
$ws.Range("M73").Value = -2564
$ws.Range("H80").Value = 2249.7144
$ws.Range("I80").Value = 1649.6
$ws.Range("K80").Value = 1649.6
$ws.Range("M80").Value = -651.5999999999999
$ws.Range("H83").Value = 2249.7144
$ws.Range("I83").Value = 1649.6
$ws.Range("K83").Value = 8248
$ws.Range("M83").Value = -3256
$ws.Range("H132").Value = 2387.889
$ws.Range("I132").Value = 1782.4286
$ws.Range("K132").Value = 5347.2858
$ws.Range("M132").Value = -2817.2858
$ws.Range("H133").Value = 110000
$ws.Range("J133").Value = 110000
$ws.Range("L133").Value = 110000
$ws.Range("N133").Value = -120120

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 18399
$ws.Range("I7").Value = 18598.285
$ws.Range("K7").Value = 18598.285
$ws.Range("M7").Value = -18486.285
$ws.Range("H93").Value = 1639.4286
$ws.Range("I93").Value = 1329.3334
$ws.Range("K93").Value = 1329.3334
$ws.Range("M93").Value = -81.33339999999998
$ws.Range("H126").Value = 18399
$ws.Range("I126").Value = 18598.285
$ws.Range("K126").Value = 55794.855
$ws.Range("M126").Value = -53324.855
$ws.Range("H136").Value = 4499.3335
$ws.Range("I136").Value = 4021.1365
$ws.Range("K136").Value = 12063.4095
$ws.Range("M136").Value = -9513.4095

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 17900
$ws.Range("I26").Value = 800
$ws.Range("J26").Value = 35000
$ws.Range("K26").Value = 800
$ws.Range("L26").Value = 35000
$ws.Range("M26").Value = -507
$ws.Range("N26").Value = -35586
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("H62").Value = 4664.923
$ws.Range("I62").Value = 3586.5
$ws.Range("K62").Value = 3586.5
$ws.Range("M62").Value = -2962.5
$ws.Range("H65").Value = 4664.923
$ws.Range("I65").Value = 3586.5
$ws.Range("K65").Value = 17932.5
$ws.Range("M65").Value = -14812.5
$ws.Range("H96").Value = 32967.5
$ws.Range("I96").Value = 20747.5
$ws.Range("K96").Value = 20747.5
$ws.Range("M96").Value = -19374.5
$ws.Range("H107").Value = 460.36365
$ws.Range("I107").Value = 280.83334
$ws.Range("K107").Value = 842.5000200000001
$ws.Range("M107").Value = 1077.49998
$ws.Range("H136").Value = 21111.908
$ws.Range("J136").Value = 54327.7
$ws.Range("L136").Value = 162983.1
$ws.Range("N136").Value = -168083.1

# --- WVR: M29 had a cached profit value; the refreshed snapshot has none ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M29").ClearContents()
